## Phase 2 preset sheet: add the "Presets Not Felt" table as its own sheet
## and populate the ranking table with participant responses.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Remove the old "Presets Not Felt" table that lived at the bottom of
#    Sheet1 (A33:B45) -- it is being moved onto its own sheet below.
# ---------------------------------------------------------------------
$oldTable = $ws1.ListObjects.Item(2)
$oldTable.Delete()
$ws1.Rows("32:34").Delete()

# ---------------------------------------------------------------------
# 2. Populate the preference ranking table (Table1) with the responses
#    collected during phase 2 testing.
# ---------------------------------------------------------------------
$rankings = @(
  @(3, 2, 4, 5, $null),
  @(2, 3, 4, 5, $null),
  @(3, 2, 4, 5, $null),
  @(4, 3, 5, 2, $null),
  @(2, 3, 1, 4, 5),
  @(3, 2, 4, 5, $null),
  @(5, 4, $null, $null, $null),
  @(3, 4, 5, $null, $null),
  @(4, 3, 5, $null, $null),
  @(3, 2, 4, 5, $null),
  @(3, 2, 4, 5, $null),
  @(4, 3, 5, $null, $null)
)

for ($i = 0; $i -lt $rankings.Length; $i++) {
  $r = 2 + $i
  $vals = $rankings[$i]
  for ($j = 0; $j -lt $vals.Length; $j++) {
    $v = $vals[$j]
    if ($v -ne $null) {
      $ws1.Cells.Item($r, 2 + $j).Value = $v
    }
  }
}

# ---------------------------------------------------------------------
# 3. Add a new worksheet ("Sheet2") after Sheet1 to hold the
#    "Presets Not Felt" responses.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)

$ws2.Range("A1").Value = "User"
$ws2.Range("B1").Value = "Presets Not Felt"

$users = @(1, 2, 3, 4, 6, 7, 8, 9, 10, 11, 12)
$notFelt = @(1, 1, 1, 1, 1, "1,2,3", "1,2", "1,2", 1, 1, "1,2")

for ($i = 0; $i -lt $users.Length; $i++) {
  $r = 2 + $i
  $ws2.Cells.Item($r, 1).Value = $users[$i]
  $ws2.Cells.Item($r, 2).Value = $notFelt[$i]
}

# ---------------------------------------------------------------------
# 4. Turn the new range into a table named "Table3" matching the one
#    that used to live on Sheet1, and right-align its data column.
# ---------------------------------------------------------------------
$newTable = $ws2.ListObjects.Add(1, $ws2.Range("A1:B1"), $null, 1)
$newTable.Name = "Table3"
$newTable.TableStyle = "TableStyleMedium9"
$newTable.Resize($ws2.Range("A1:B12"))

$ws2.Range("B2:B12").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 5. Restore the selections on each sheet, Sheet2 ends up active.
# ---------------------------------------------------------------------
$ws1.Range("J27").Select()
$ws2.Activate()
$ws2.Range("A1:XFD1").Select()
